$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.635.83"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3
$ws.Range("D3").Value = "3.692.62"
$ws.Range("E3").Value = "  +0.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'672.62"
$ws.Range("E5").Value = "  -1.35%  "

# Row 6
$ws.Range("D6").Value = "'161.24"
$ws.Range("E6").Value = "  +1.75%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.499"
$ws.Range("E8").Value = "  +1.08%  "

# Row 9
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$ws.Range("D10").Value = "'7.12"
$ws.Range("E10").Value = "  +1.81%  "

# Row 11
$ws.Range("D11").Value = "'0.443"
$ws.Range("E11").Value = "  +1.21%  "

# Row 12
$ws.Range("E12").Value = "  +0.92%  "

# Row 13
$ws.Range("D13").Value = "'33.40"
$ws.Range("E13").Value = "  +2.98%  "

# Row 14
$ws.Range("D14").Value = "3.668.95"
$ws.Range("E14").Value = "  +0.20%  "

# Row 15
$ws.Range("D15").Value = "69.615.87"
$ws.Range("E15").Value = "  +0.50%  "

# Row 16
$ws.Range("D16").Value = "'0.117"
$ws.Range("E16").Value = "  +1.55%  "

# Row 17
$ws.Range("D17").Value = "'16.20"
$ws.Range("E17").Value = "  +0.94%  "

# Row 18
$ws.Range("D18").Value = "'6.51"
$ws.Range("E18").Value = "  +1.19%  "

# Row 19
$ws.Range("D19").Value = "'472.91"
$ws.Range("E19").Value = "  +0.50%  "

# Row 20
$ws.Range("D20").Value = "'9.80"
$ws.Range("E20").Value = "  -1.75%  "

# Row 21
$ws.Range("D21").Value = "'0.651"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
$ws.Range("D22").Value = "'80.18"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
$ws.Range("D23").Value = "3.840.83"
$ws.Range("E23").Value = "  +0.35%  "

# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "'0.0000128"
$ws.Range("E24").Value = "  +4.16%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("D26").Value = "'11.01"
$ws.Range("E26").Value = "  +0.38%  "

# Row 27
$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "  +0.46%  "

# Row 28
$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$ws.Range("D29").Value = "'1.73"
$ws.Range("E29").Value = "  -0.95%  "

# Row 30
$ws.Range("D30").Value = "'2.03"
$ws.Range("E30").Value = "  +2.11%  "

# Row 31
$ws.Range("D31").Value = "'0.170"
$ws.Range("E31").Value = "  +4.24%  "

# Row 32
$ws.Range("E32").Value = "  +0.35%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.54"
$ws.Range("E33").Value = "  -1.32%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.91"
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("D35").Value = "3.687.98"
$ws.Range("E35").Value = "  +0.66%  "

# Row 36
$ws.Range("D36").Value = "'8.52"
$ws.Range("E36").Value = "  +3.75%  "

# Row 37
$ws.Range("E37").Value = "  -0.75%  "

# Row 38
$ws.Range("E38").Value = "  -0.01%  "

# Row 39
$ws.Range("D39").Value = "'2.28"
$ws.Range("E39").Value = "  +1.19%  "

# Row 40
$ws.Range("E40").Value = "  -0.02%  "

# Row 41
$ws.Range("D41").Value = "'0.0913"
$ws.Range("E41").Value = "  +1.06%  "

# Row 42
$ws.Range("D42").Value = "'175.77"
$ws.Range("E42").Value = "  +1.22%  "

# Row 43
$ws.Range("E43").Value = "  -0.50%  "

# Row 44
$ws.Range("D44").Value = "'47.00"
$ws.Range("E44").Value = "  -1.13%  "

# Row 45
$ws.Range("D45").Value = "'2.77"
$ws.Range("E45").Value = "  +2.49%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'28.03"
$ws.Range("E46").Value = "  +1.80%  "

# Row 47
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "'0.000276"
$ws.Range("E47").Value = "  -2.20%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.29"
$ws.Range("E48").Value = "  +0.89%  "

# Row 49
$ws.Range("D49").Value = "'1.09"
$ws.Range("E49").Value = "  -0.13%  "

# Row 50
$ws.Range("D50").Value = "'7.90"
$ws.Range("E50").Value = "  +1.29%  "

# Row 51
$ws.Range("E51").Value = "  -0.28%  "
